$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("testSheet1")

# Move the data currently in rows 3:4 down to rows 5:6, leaving rows 3:4 empty.
$ws1.Range("A3:I4").Copy($ws1.Range("A5"))
$ws1.Range("K3:L4").Copy($ws1.Range("K5"))
$ws1.Range("N3:N4").Copy($ws1.Range("N5"))
$ws1.Range("A3:N4").Clear()

# Activate testSheet1 (becomes the tab shown when the workbook is opened)
# and select cell B4 on it, matching the new selection state.
$ws1.Activate() | Out-Null
$ws1.Range("B4").Select() | Out-Null
